$wb = $excel.ActiveWorkbook

# 展览 (Exhibition) sheet updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 20264
$wsExpo.Range("F4").Value = 803
$wsExpo.Range("F5").Value = 320
$wsExpo.Range("F6").Value = 1101
$wsExpo.Range("F7").Value = 17
$wsExpo.Range("F8").Value = 7629
$wsExpo.Range("F9").Value = 522
$wsExpo.Range("F10").Value = 737
$wsExpo.Range("F11").Value = 273
$wsExpo.Range("F13").Value = 162
$wsExpo.Range("F14").Value = 128
$wsExpo.Range("F15").Value = 15
$wsExpo.Range("F17").Value = 196
$wsExpo.Range("F18").Value = 1342
$wsExpo.Range("F19").Value = 440
$wsExpo.Range("F21").Value = 688
$wsExpo.Range("F24").Value = 70
$wsExpo.Range("F25").Value = 326
$wsExpo.Range("F26").Value = 1119
$wsExpo.Range("F27").Value = 35
$wsExpo.Range("F30").Value = 5218
$wsExpo.Range("F31").Value = 567
$wsExpo.Range("F32").Value = 80
$wsExpo.Range("F33").Value = 2872
$wsExpo.Range("F37").Value = 12693
$wsExpo.Range("F38").Value = 1340
$wsExpo.Range("F39").Value = 90
$wsExpo.Range("F41").Value = 60
$wsExpo.Range("F42").Value = 274
$wsExpo.Range("F43").Value = 380
$wsExpo.Range("F44").Value = 4008
$wsExpo.Range("F45").Value = 322

# 演出 (Performance) sheet updates
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 199

# 全部类型 (All types) sheet updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 20266
$wsAll.Range("F4").Value = 803
$wsAll.Range("F5").Value = 320
$wsAll.Range("F6").Value = 1101
$wsAll.Range("F7").Value = 17
$wsAll.Range("F8").Value = 7629
$wsAll.Range("F9").Value = 522
$wsAll.Range("F10").Value = 737
$wsAll.Range("F11").Value = 274
$wsAll.Range("F13").Value = 162
$wsAll.Range("F14").Value = 128
$wsAll.Range("F15").Value = 15
$wsAll.Range("F17").Value = 196
$wsAll.Range("F18").Value = 1342
$wsAll.Range("F19").Value = 440
$wsAll.Range("F21").Value = 688
$wsAll.Range("F24").Value = 70
$wsAll.Range("F25").Value = 326
$wsAll.Range("F26").Value = 1119
$wsAll.Range("F27").Value = 35
$wsAll.Range("F30").Value = 199
$wsAll.Range("F31").Value = 5218
$wsAll.Range("F32").Value = 567
$wsAll.Range("F34").Value = 80
$wsAll.Range("F36").Value = 2872
$wsAll.Range("F40").Value = 12693
$wsAll.Range("F41").Value = 1340
$wsAll.Range("F42").Value = 90
$wsAll.Range("F44").Value = 60
$wsAll.Range("F45").Value = 274
$wsAll.Range("F46").Value = 380
$wsAll.Range("F47").Value = 4008
$wsAll.Range("F48").Value = 322
